$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.646.54"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.327.39"
$ws.Range("E3").Value = "  +4.36%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "270.85"
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.16"
$ws.Range("E6").Value = "  +8.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  +2.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.72"
$ws.Range("E10").Value = "  -1.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +2.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.01"
$ws.Range("E12").Value = "  +4.34%  "

$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.677.24"
$ws.Range("E14").Value = "  +4.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.72"
$ws.Range("E15").Value = "  +5.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.856"
$ws.Range("E16").Value = "  +8.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.335.69"
$ws.Range("E17").Value = "  +3.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.623.96"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("E19").Value = "  +3.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.31"
$ws.Range("E20").Value = "  +5.67%  "

$ws.Range("E21").Value = "  +2.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.85"
$ws.Range("E22").Value = "  +2.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.25"
$ws.Range("E23").Value = "  -4.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.57"
$ws.Range("E24").Value = "  +9.37%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("E26").Value = "  -1.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  +4.40%  "

$ws.Range("E28").Value = "  -1.80%  "

$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.56"
$ws.Range("E30").Value = "  -1.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.52"
$ws.Range("E31").Value = "  +8.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "172.71"
$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0894"
$ws.Range("E33").Value = "  -0.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.47"
$ws.Range("E34").Value = "  +2.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.127"
$ws.Range("E35").Value = "  +2.91%  "

$ws.Range("E36").Value = "  +0.61%  "

$ws.Range("E37").Value = "  -3.60%  "

$ws.Range("E38").Value = "  +1.83%  "

$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.38"
$ws.Range("E40").Value = "  +9.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.233"
$ws.Range("E41").Value = "  +11.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.36"
$ws.Range("E42").Value = "  +20.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.02"
$ws.Range("E43").Value = "  -3.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.08"
$ws.Range("E44").Value = "  +6.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.59"
$ws.Range("E45").Value = "  -3.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.36"
$ws.Range("E46").Value = "  -0.73%  "

$ws.Range("E47").Value = "  +5.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.76"
$ws.Range("E48").Value = "  +0.31%  "

$ws.Range("E49").Value = "  +1.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.554.33"
$ws.Range("E50").Value = "  +4.36%  "

$ws.Range("E51").Value = "  +12.57%  "
